# cv126132a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had two label-only rows ("situação do domicílio" at row 5 and
# "grandes regiões" at row 8) that separated groups of data rows. The fix
# removes those two spacer/header rows (letting the data rows that followed
# them shift up into their place) and renames the "unnamed: 1_level_1"
# column header (B2) to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("situação do domicílio") was a label-only row with no data - delete it
# so the "urbana" row (previously row 6) moves up to row 5.
$ws.Rows.Item(5).Delete()

# After the row-5 deletion, the "grandes regiões" label-only row (previously
# row 8) is now row 7 - delete it too, so "norte" (previously row 9) moves
# up to row 7.
$ws.Rows.Item(7).Delete()

# Rename the "unnamed: 1_level_1" sub-header to "total".
$ws.Range("B2").Value = "total"
